# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp note (A1) from 21:22 to 21:52
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 21:52"

# A new "Asturias" entry is inserted (keeping the list sorted by Casos totales,
# descending). This pushes "Gipuzkoa/Guipuzcoa" (formerly row 23) down to row 24
# with its data unchanged, while row 23 now carries Asturias with fresh figures.
$ws.Range("A23").Value = "Asturias"
$ws.Range("B23").Value = 1799
$ws.Range("C23").Value = 372
$ws.Range("D23").Value = 1306
$ws.Range("E23").Value = 121

$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B24").Value = 1756
$ws.Range("C24").Value = 4514
$ws.Range("D24").Value = 4603
$ws.Range("E24").Value = 103

# Melilla (row 54) figures updated
$ws.Range("B54").Value = 95
$ws.Range("C54").Value = 17
$ws.Range("D54").Value = 76
